# Update "想去人数" (want-to-go count) values in column F across all sheets
# to match the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1784
$ws1.Range("F4").Value  = 452
$ws1.Range("F7").Value  = 631
$ws1.Range("F8").Value  = 333
$ws1.Range("F9").Value  = 1728
$ws1.Range("F10").Value = 365
$ws1.Range("F11").Value = 1423
$ws1.Range("F13").Value = 336
$ws1.Range("F14").Value = 678
$ws1.Range("F15").Value = 12795
$ws1.Range("F16").Value = 12797
$ws1.Range("F17").Value = 954
$ws1.Range("F22").Value = 563
$ws1.Range("F23").Value = 2004
$ws1.Range("F24").Value = 29
$ws1.Range("F27").Value = 19
$ws1.Range("F28").Value = 247
$ws1.Range("F29").Value = 677

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 81
$ws2.Range("F6").Value  = 16
$ws2.Range("F10").Value = 77

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 86

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 86
$ws4.Range("F5").Value  = 1784
$ws4.Range("F6").Value  = 452
$ws4.Range("F11").Value = 631
$ws4.Range("F13").Value = 333
$ws4.Range("F14").Value = 1728
$ws4.Range("F15").Value = 365
$ws4.Range("F16").Value = 1423
$ws4.Range("F18").Value = 336
$ws4.Range("F19").Value = 81
$ws4.Range("F20").Value = 678
$ws4.Range("F21").Value = 12795
$ws4.Range("F22").Value = 12797
$ws4.Range("F23").Value = 954
$ws4.Range("F28").Value = 563
$ws4.Range("F29").Value = 16
$ws4.Range("F31").Value = 2004
$ws4.Range("F32").Value = 29
$ws4.Range("F37").Value = 19
$ws4.Range("F38").Value = 247
$ws4.Range("F39").Value = 677
$ws4.Range("F40").Value = 77

$wb.Save()
